$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Property" value for the Loads row (B3): NumPhases -> VoltagesMagAng
$ws.Range("B3").Value = "VoltagesMagAng"

# Remove the third table column ("Column1"/kW data) from Table1, then delete
# the now-empty worksheet column so the grid collapses back down.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item("Column1").Delete()
$ws.Columns.Item(3).Delete()

# Move the active selection, matching the saved view state.
[void]$ws.Range("B9").Select()
